$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("заказы")

$data = @(
    @(26, 1, 4, 6, 2),
    @(27, 2, 1464, 1, 12),
    @(28, 3, 12221, 1, 11)
)

$startRow = 8
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $values = $data[$i]
    $ws.Cells.Item($row, 1).Value = $values[0]
    $ws.Cells.Item($row, 2).Value = $values[1]
    $ws.Cells.Item($row, 3).Value = $values[2]
    $ws.Cells.Item($row, 4).Value = $values[3]
    $ws.Cells.Item($row, 5).Value = $values[4]
}
